$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General-looking numeric strings must stay literal text, matching source data)
$textCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D26", "D27", "D29", "D32", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.132.70"
$ws.Range("D3").Value = "2.982.21"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "576.86"
$ws.Range("E5").Value = "  -0.83%  "
$ws.Range("D6").Value = "162.18"
$ws.Range("E6").Value = "  +6.30%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").Value = "2.977.81"
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "6.70"
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").Value = "34.64"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "66.109.61"
$ws.Range("E16").Value = "  +2.77%  "
$ws.Range("D17").Value = "3.474.85"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").Value = "6.91"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").Value = "2.981.74"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("D20").Value = "450.83"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "13.84"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").Value = "7.28"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").Value = "82.03"
$ws.Range("E24").Value = "  +1.91%  "
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "12.22"
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "10.10"
$ws.Range("E27").Value = "  -7.44%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "8.15"
$ws.Range("E29").Value = "  +6.94%  "
$ws.Range("E30").Value = "  +12.39%  "
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").Value = "0.0000101"
$ws.Range("E32").Value = "  -6.79%  "
$ws.Range("D33").Value = "27.24"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("D36").Value = "0.983"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  +2.95%  "
$ws.Range("D38").Value = "2.05"
$ws.Range("E38").Value = "  -4.38%  "
$ws.Range("D39").Value = "49.46"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("D40").Value = "44.16"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("D41").Value = "2.87"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("E43").Value = "  +0.74%  "
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("D45").Value = "389.08"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").Value = "0.0356"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D47").Value = "2.729.51"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "131.64"
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "23.30"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.106"
$ws.Range("E51").Value = "  +0.89%  "

# Restore default style (no explicit number format) now that the text values are locked in
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
